$p = $ppt.ActivePresentation

# ---------------------------------------------------------------------
# 1) Update the "Insert Header & Footer" fixed date / slide-number text
#    on the slide master and every slide layout:
#      date placeholder  : "3/3/2017" -> "30/05/17"
#      slide# placeholder: "<N deg>"  -> "<#>"   (angle-quote glyphs)
# ---------------------------------------------------------------------
function Update-HeaderFooterShapes($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        $phType = $null
        try { $phType = $sh.PlaceholderFormat.Type } catch {}
        if ($phType -eq 16) {
            # ppPlaceholderDate
            $sh.TextFrame.TextRange.Text = "30/05/17"
        } elseif ($phType -eq 13) {
            # ppPlaceholderSlideNumber
            $sh.TextFrame.TextRange.Text = [char]0x2039 + "#" + [char]0x203A
        }
    }
}

Update-HeaderFooterShapes($p.SlideMaster.Shapes)
for ($L = 1; $L -le $p.SlideMaster.CustomLayouts.Count; $L++) {
    Update-HeaderFooterShapes($p.SlideMaster.CustomLayouts.Item($L).Shapes)
}

# ---------------------------------------------------------------------
# 2) Slide 1 subtitle: the session date "6 mars 2017" -> "31 mai 2017"
#    (edited in two passes, day then month, to mirror the original
#    author's edit and its resulting run split)
# ---------------------------------------------------------------------
$s1 = $p.Slides.Item(1)

$subtitle = $null
for ($i = 1; $i -le $s1.Shapes.Count; $i++) {
    $sh = $s1.Shapes.Item($i)
    $phType = $null
    try { $phType = $sh.PlaceholderFormat.Type } catch {}
    if ($phType -eq 4) {
        # ppPlaceholderSubtitle
        $subtitle = $sh
    }
}

$tr = $subtitle.TextFrame.TextRange
$dateParaIdx = -1
for ($i = 1; $i -le $tr.Paragraphs().Count; $i++) {
    $para = $tr.Paragraphs($i, 1)
    if ($para.Text -like "*mars*") {
        $dateParaIdx = $i
    }
}

# "6 mars 2017" -> replace day "6" with "31"
$dateParagraph = $tr.Paragraphs($dateParaIdx, 1)
$day = $dateParagraph.Characters(1, 1)
$day.Text = "31"

# re-fetch paragraph (length changed) -> replace " mars " with " mai "
$dateParagraph2 = $tr.Paragraphs($dateParaIdx, 1)
$monthStart = $dateParagraph2.Text.IndexOf(" mars ") + 1
$month = $dateParagraph2.Characters($monthStart, 6)
$month.Text = " mai "
